$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Invalid and Absent flags set
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4: Total Attendance Count and Real flags set
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5: Total Attendance Count and Real flags set
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6: Absent flag set
$ws.Range("H6").Value = 1

# Row 7: Absent flag set
$ws.Range("H7").Value = 1

# Row 8: Absent flag set
$ws.Range("H8").Value = 1

# Row 9: Total Attendance Count and Real flags set
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10: Absent flag set
$ws.Range("H10").Value = 1

# Row 11: Absent flag set
$ws.Range("H11").Value = 1

# Row 12: Absent flag set
$ws.Range("H12").Value = 1

# Row 13: Absent flag set
$ws.Range("H13").Value = 1

# Row 14: Absent flag set
$ws.Range("H14").Value = 1

# Row 15: Absent flag set
$ws.Range("H15").Value = 1

# Row 16: Absent flag set
$ws.Range("H16").Value = 1

# Row 17: Absent flag set
$ws.Range("H17").Value = 1

# Row 18: Absent flag set
$ws.Range("H18").Value = 1
